# Weekly update: a new price observation is inserted at the top of the
# "Orégano" / "Vega Central Mapocho de Santiago" data block (row 67),
# pushing the existing rows 67-95 down to 68-96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 67; Excel shifts rows 67:95 down to 68:96
# and extends the used range/dimension automatically.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with this week's data.
$ws.Range("A67").Value = 9
$ws.Range("B67").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C67").Value = "Metropolitana"
$ws.Range("D67").Value = 45007
$ws.Range("E67").Value = 13
$ws.Range("F67").Value = 100112029
$ws.Range("G67").Value = "Orégano"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 16
$ws.Range("K67").Value = 17000
$ws.Range("L67").Value = 18000
$ws.Range("M67").Value = 17500
$ws.Range("N67").Value = "$/docena de atados"
$ws.Range("O67").Value = "Región Metropolitana"
$ws.Range("P67").Value = 5833
$ws.Range("Q67").Value = 3
$ws.Range("R67").Value = "Hortaliza"
